$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.080.37"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +0.01%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.816.98"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -0.61%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9986"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.22%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "311.35"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +0.08%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9987"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -0.12%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4431"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +4.76%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3737"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +1.66%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07440"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +2.88%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.8723"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +3.41%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "20.83"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +0.35%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.798.64"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -1.55%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "6.696"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +0.43%  "
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +3.90%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.07097"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +0.82%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "5.320"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +0.44%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.9987"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000008736"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -0.25%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.9990"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -0.21%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "14.95"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +0.38%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "27.097.52"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +0.00%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.210"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +1.31%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "10.90"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +0.73%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.026.55"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -0.99%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.977"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -0.44%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "151.56"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -0.12%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.242"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -0.07%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "18.51"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +1.18%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.276"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +0.21%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "118.13"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +1.04%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.08822"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +1.16%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.7606"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +2.96%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.176"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -0.39%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.541"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +2.55%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.885"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -0.44%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.9980"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -0.16%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.100"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +0.71%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01981"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +1.92%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.05261"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +0.12%  "
$ws.Range("B40").Value = "FraxShare"
$ws.Range("C40").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "7.208"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -1.66%  "
$ws.Range("B41").Value = "TheSandbox"
$ws.Range("C41").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.5284"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +4.16%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.841"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -1.32%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.1708"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +1.32%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.180"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +11.43%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "8.641"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +1.01%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.5024"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +6.18%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "10.53"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +0.08%  "
$ws.Range("B48").Value = "NEARProtocol"
$ws.Range("C48").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.693"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +2.19%  "
$ws.Range("B49").Value = "Quant"
$ws.Range("C49").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "105.04"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -1.07%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.9981"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -0.11%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.06353"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +0.28%  "
